$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of data (row 21) after the last existing data row (row 20)
# Force column A to be plain text so the date string isn't auto-converted
# into a date serial number, matching the other "MM/DD/YYYY" text cells.
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "01/02/2026"
$ws.Cells.Item(21, 1).Style = "Normal"

$ws.Cells.Item(21, 2).Value = 1101.564
$ws.Cells.Item(21, 3).Value = 0.04493610902317068
$ws.Cells.Item(21, 4).Value = 50
